$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL_Test_Cases")

# Clear the "Status" column (Not Executed) for the data rows.
$ws.Range("K4:K8").ClearContents() | Out-Null

# Update the selection / scroll position to match the saved view (A8 selected, scrolled to top-left).
$ws.Activate() | Out-Null
$ws.Range("A8").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
